# Apply the "Long term project planning" sprint-objectives edits.
$wb = $excel.ActiveWorkbook

# --- Sheet: "Comms_Notification prefs" ----------------------------------
$wsComms = $wb.Worksheets.Item("Comms_Notification prefs")
$wsComms.Select()
$wsComms.Range("B9").Select()

# --- Sheet: "General auth exp work" ------------------------------------
# Rows 4-6 content (height + value) shifts down to rows 5-7; row 4 becomes "N/A".
$wsAuth = $wb.Worksheets.Item("General auth exp work")

$valB4 = $wsAuth.Range("B4").Value()
$valB5 = $wsAuth.Range("B5").Value()
$valB6 = $wsAuth.Range("B6").Value()

$wsAuth.Range("B7").Value = $valB6
$wsAuth.Range("B7").WrapText = $true
$wsAuth.Rows.Item(7).RowHeight = 34

$wsAuth.Range("B6").Value = $valB5
$wsAuth.Rows.Item(6).RowHeight = 51

$wsAuth.Range("B5").Value = $valB4
$wsAuth.Rows.Item(5).RowHeight = 34

$wsAuth.Range("A4").Copy()
$wsAuth.Range("B4").PasteSpecial(-4122)
$wsAuth.Range("B4").Value = "N/A"
$wsAuth.Rows.Item(4).AutoFit()

$wsAuth.Select()
$wsAuth.Range("B3:B4").Select()

# --- Sheet: "Long term resourcing by person" ---------------------------
$wsRes = $wb.Worksheets.Item("Long term resourcing by person")
$wsRes.Range("D7").Value = 0.4
$wsRes.Range("C9").Value = 0
$wsRes.Range("F9").Value = 0.2
$wsRes.Range("B10").Value = 0.15

# C10 was a "TBD" text cell; it becomes a percentage value like its neighbors.
$wsRes.Range("B10").Copy()
$wsRes.Range("C10").PasteSpecial(-4122)
$wsRes.Range("C10").Value = 0.15

$wsRes.Select()
$wsRes.Range("D9").Select()
